$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.622.34'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = '3.698.09'
$ws.Range('E3').Value = '  +0.54%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''670.76'
$ws.Range('E5').Value = '  -2.00%  '

$ws.Range('D6').Value = '''161.40'
$ws.Range('E6').Value = '  +1.52%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('E8').Value = '  +0.80%  '

$ws.Range('E9').Value = '  +0.39%  '

$ws.Range('D10').Value = '''7.07'
$ws.Range('E10').Value = '  +0.58%  '

$ws.Range('E11').Value = '  +1.96%  '

$ws.Range('E12').Value = '  +0.94%  '

$ws.Range('D13').Value = '''32.79'
$ws.Range('E13').Value = '  +1.61%  '

$ws.Range('D14').Value = '3.650.35'
$ws.Range('E14').Value = '  -0.67%  '

$ws.Range('D15').Value = '69.671.38'

$ws.Range('E16').Value = '  +1.54%  '

$ws.Range('E17').Value = '  +2.53%  '

$ws.Range('E18').Value = '  +1.71%  '

$ws.Range('D19').Value = '''474.18'
$ws.Range('E19').Value = '  +1.00%  '

$ws.Range('D20').Value = '''9.80'
$ws.Range('E20').Value = '  -1.43%  '

$ws.Range('D21').Value = '''0.653'
$ws.Range('E21').Value = '  +0.65%  '

$ws.Range('D22').Value = '''80.41'
$ws.Range('E22').Value = '  +0.69%  '

$ws.Range('D23').Value = '3.846.96'
$ws.Range('E23').Value = '  +0.59%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '''0.0000127'
$ws.Range('E25').Value = '  +3.53%  '

$ws.Range('D26').Value = '''10.94'
$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('D27').Value = '''9.11'
$ws.Range('E27').Value = '  -0.98%  '

$ws.Range('D28').Value = '''2.68'
$ws.Range('E28').Value = '  -0.83%  '

$ws.Range('D29').Value = '''1.74'
$ws.Range('E29').Value = '  +0.61%  '

$ws.Range('E30').Value = '  +1.54%  '

$ws.Range('D31').Value = '''6.61'
$ws.Range('E31').Value = '  +0.99%  '

$ws.Range('E32').Value = '  +4.35%  '

$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  -0.11%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '''26.89'
$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('D35').Value = '3.686.15'
$ws.Range('E35').Value = '  +0.89%  '

$ws.Range('D36').Value = '''8.55'
$ws.Range('E36').Value = '  +5.22%  '

$ws.Range('D37').Value = '''6.08'
$ws.Range('E37').Value = '  -1.01%  '

$ws.Range('E39').Value = '  +0.78%  '

$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('D41').Value = '''0.0909'
$ws.Range('E41').Value = '  +1.15%  '

$ws.Range('D42').Value = '''172.26'
$ws.Range('E42').Value = '  +3.62%  '

$ws.Range('D43').Value = '''0.939'
$ws.Range('E43').Value = '  +0.01%  '

$ws.Range('D44').Value = '''47.05'
$ws.Range('E44').Value = '  -0.94%  '

$ws.Range('D45').Value = '''2.78'
$ws.Range('E45').Value = '  +2.76%  '

$ws.Range('D46').Value = '''0.000280'
$ws.Range('E46').Value = '  -0.66%  '

$ws.Range('D47').Value = '''27.69'
$ws.Range('E47').Value = '  +1.87%  '

$ws.Range('E48').Value = '  -1.60%  '

$ws.Range('E49').Value = '  -1.18%  '

$ws.Range('D50').Value = '''7.89'
$ws.Range('E50').Value = '  +1.31%  '

$ws.Range('E51').Value = '  +0.58%  '
